$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.657.52'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '2.928.42'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = '2.928.93'
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.500'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.141'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.441'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000221'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.127'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = '3.410.81'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').Value = '61.487.92'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '2.976.17'
$ws.Range('E18').Value = '  +1.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '434.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.667'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.09%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -8.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.20%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.05'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('D35').Value = '0.0₃0853'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.981'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.20'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.14%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.46'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.117'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.272'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.58%  '
$ws.Range('D45').Value = '2.681.48'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0336'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.94%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '132.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '357.80'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.97%  '
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.52%  '
